# Revert "Revert "typo: Database 2 -> Database 3""
#
# The deck still has a handful of slide titles reading "Database 2: ..."
# that should say "Database 3: ..." instead. Walk every slide/shape,
# find the title run(s) that start with "Database 2:" and fix them up,
# rewriting only the minimal amount of text so existing run formatting
# (and any later runs in the same paragraph) stay untouched.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text
        if (-not $fullText) { continue }
        if ($fullText.IndexOf("Database 2:") -lt 0) { continue }

        # There is exactly one "Database 2: " occurrence per affected
        # title, always at the very start of the text box.
        if ($fullText.StartsWith("Database 2: ")) {
            $prefix = "Database 2: "
            $firstRun = $tr.Characters(1, $prefix.Length)

            if ($firstRun.Text -eq $prefix) {
                # The run that currently holds "Database 2: " continues
                # on with more text (e.g. "ucddb002 128Hz original (+ ",
                # or just "100Hz - all ..."). If the WHOLE run is exactly
                # the prefix, split "Database 3: " into its own run and
                # drop the stale "Database 2: " text from the remainder
                # (slide 9's case). Otherwise the prefix is only part of
                # a longer first run, so just rewrite that whole run's
                # text in place (no run split, slides 2-5 & 8's case).
                $runLen = $firstRun.Length
                if ($runLen -eq $prefix.Length) {
                    [void]$firstRun.InsertBefore("Database 3: ")
                    $staleePrefix = $tr.Characters($prefix.Length + 1, $prefix.Length)
                    $staleePrefix.Text = ""
                } else {
                    $wholeRun = $tr.Characters(1, $runLen)
                    $oldRunText = $wholeRun.Text
                    if ($oldRunText.StartsWith($prefix)) {
                        $newRunText = "Database 3: " + $oldRunText.Substring($prefix.Length)
                        $wholeRun.Text = $newRunText
                    }
                }
            }
        }
    }
}
